$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.402.17"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.692.77"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5489"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.63%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2733"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06471"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07674"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.94%  "
$ws.Range("D12").Value = "1.700.18"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.553"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008394"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "26.443.84"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.951"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.011"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("E25").Value = "  +7.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.901"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("E28").Value = "  -5.10%  "
$ws.Range("E29").Value = "  +2.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.331"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.604"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6157"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.410"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.709"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.220"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D39").Value = "1.119.78"
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("E40").Value = "  +0.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8860"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("D44").Value = "1.843.19"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("E46").Value = "  -5.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.228"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05282"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.104"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.83%  "
$ws.Range("E51").Value = "  +0.09%  "
